# Update "想去人数" (interest count) figures in the "展览" and "全部类型"
# sheets to match the freshly generated gh-pages output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1147
$ws1.Range("F5").Value = 187
$ws1.Range("F6").Value = 10
$ws1.Range("F8").Value = 266
$ws1.Range("F9").Value = 391
$ws1.Range("F13").Value = 547
$ws1.Range("F14").Value = 162
$ws1.Range("F15").Value = 12996
$ws1.Range("F16").Value = 153
$ws1.Range("F17").Value = 13
$ws1.Range("F18").Value = 5337
$ws1.Range("F19").Value = 5543

# --- Sheet "全部类型" (all categories) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1147
$ws4.Range("F5").Value = 187
$ws4.Range("F6").Value = 10
$ws4.Range("F9").Value = 266
$ws4.Range("F10").Value = 391
$ws4.Range("F15").Value = 547
$ws4.Range("F16").Value = 162
$ws4.Range("F17").Value = 12996
$ws4.Range("F18").Value = 153
$ws4.Range("F20").Value = 13
$ws4.Range("F21").Value = 5337
$ws4.Range("F22").Value = 5543
